# "Mar Time sheet added" -> actually resets the visible "Time sheet " tab
# from a March-2024 period back to a February-2024 period (leap year,
# 1-Feb-2024 .. 29-Feb-2024), re-derives the weekday labels, and
# recomputes the per-day Status / Shift Timing / Hours-worked columns
# from the (now corrected) "Date of Joining" (D5) and weekend pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$wf = $excel.WorksheetFunction

# ---- header dates -------------------------------------------------------
# F1 : "Timesheet for the Month" banner date
# D5 : Date of Joining
# B9 : period From, C9 : period To
$ws.Range("F1").Value = 45323
$ws.Range("D5").Value = 45336
$ws.Range("B9").Value = 45323
$ws.Range("C9").Value = 45351

$joinDate = 45336

# ---- daily grid (rows 11..39, A=Sr.No, B=Date formula, C=Day formula) ---
# B11:B38 / C11:C39 are (mostly) shared formulas that recompute
# automatically once B9/C9 change. Row 39 (the trailing 29th day) is the
# one exception: B39 is a hard-coded date and C39 used to be a literal
# string value ("Friday") instead of the =TEXT(B39,"dddd") formula used
# by every other row in the block - fix that up so it behaves the same
# as C11:C38.
$ws.Range("C39").Formula = '=TEXT(B39,"dddd")'

$null = $wb.Application.Calculate()

for ($row = 11; $row -le 39; $row++) {
    $dateCell = $ws.Cells.Item($row, 2)
    $serial = $dateCell.Value2()

    $dCell = $ws.Cells.Item($row, 4)
    $eCell = $ws.Cells.Item($row, 5)
    $fCell = $ws.Cells.Item($row, 6)

    if ($serial -eq "") {
        continue
    }

    if ($serial -lt $joinDate) {
        $dCell.Value = "NA"
        $null = $eCell.ClearContents()
        $null = $fCell.ClearContents()
    }
    else {
        $weekday = $wf.Weekday($serial, 2)
        if ($weekday -ge 6) {
            $dCell.Value = "WeekOff"
            $null = $eCell.ClearContents()
            $null = $fCell.ClearContents()
        }
        else {
            $dCell.Value = "Present"
            $eCell.Value = "00:9 AM-To-6:00PM"
            $fCell.Value = "09:00Hours"
        }
    }
}

# ---- selection, as left by the edit -------------------------------------
$null = $ws.Range("D5:F5").Select()

$null = $wb.Application.Calculate()
